# Update the cryptos price list: refresh Price/Volume(1h) figures for the
# existing rows and insert a new "LEO" entry at row 26 (pushing the rest of
# the list down by one row and dropping the final "Aave" row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.499.96"
$ws.Range("D3").Value = "1.680.09"
$ws.Range("E3").Value = "  +5.75%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9964"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3686"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3421"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.10"
$ws.Range("E9").Value = "  +16.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.160"
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07223"
$ws.Range("E11").Value = "  +4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9983"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.103"
$ws.Range("E13").Value = "  +4.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.16"
$ws.Range("E14").Value = "  +3.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.709"
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").Value = "1.677.49"
$ws.Range("E16").Value = "  +5.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").Value = "  +3.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9962"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06652"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "80.60"
$ws.Range("E20").Value = "  +6.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.42"
$ws.Range("E21").Value = "  +3.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.074"
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.09"
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("D24").Value = "24.436.11"
$ws.Range("E24").Value = "  +10.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.423"
$ws.Range("E25").Value = "  +2.31%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.350"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.652"
$ws.Range("E27").Value = "  +6.39%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.12"
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.42"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").Value = "1.863.69"
$ws.Range("E30").Value = "  +5.61%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.03"
$ws.Range("E31").Value = "  +4.63%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.241"
$ws.Range("E32").Value = "  +6.17%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.021"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9785"
$ws.Range("E34").Value = "  +6.81%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08413"
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.692"
$ws.Range("E36").Value = "  +4.74%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.33"
$ws.Range("E37").Value = "  +5.79%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06370"
$ws.Range("E38").Value = "  +5.94%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.314"
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.674"
$ws.Range("E40").Value = "  +4.26%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02305"
$ws.Range("E41").Value = "  +6.01%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.245"
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2090"
$ws.Range("E43").Value = "  +5.63%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6084"
$ws.Range("E44").Value = "  +5.52%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9965"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.760"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.94"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5860"
$ws.Range("E48").Value = "  +5.51%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.35"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.006"
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07182"
$ws.Range("E51").Value = "  +7.02%  "
